# Ran code for averaged intensities on spiral schemes.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BrassA")

# Row 2 (C2:P2) HKL-index header labels (text is unchanged, only the
# underlying shared-string bookkeeping shifts once the new scheme rows
# below are written).
$ws.Range("B2").Value = "HKL"
$ws.Range("C2").Value = "[1, 1, 1]"
$ws.Range("D2").Value = "[2, 0, 0]"
$ws.Range("E2").Value = "[2, 2, 0]"
$ws.Range("F2").Value = "[3, 1, 1]"
$ws.Range("G2").Value = "[2, 2, 2]"
$ws.Range("H2").Value = "[4, 0, 0]"
$ws.Range("I2").Value = "[3, 3, 1]"
$ws.Range("J2").Value = "[4, 2, 0]"
$ws.Range("K2").Value = "[4, 2, 2]"
$ws.Range("L2").Value = "[5, 1, 1]"
$ws.Range("M2").Value = "[3, 3, 3]"
$ws.Range("N2").Value = "2Pairs"
$ws.Range("O2").Value = "4Pairs"
$ws.Range("P2").Value = "MaxUnique"

# New averaging schemes (Gaussian-Quadrature + the 3 spiral schemes) are
# inserted into the scheme list, pushing the remaining scheme names down.
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("B16").Value = "Rotation-60detTilt"

# Three new rows (17-19) with averaged intensities (all 1) for the
# remaining schemes, appended after the previous last row.
$newRows = @(
    @{ Row = 17; A = 15; B = "HexGrid-90degTilt5degRes" },
    @{ Row = 18; A = 16; B = "HexGrid-90degTilt22p5degRes" },
    @{ Row = 19; A = 17; B = "HexGrid-60degTilt5degRes" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    $prevRow = $rowNum - 1

    # Copy formatting from the previous row so the new row matches the
    # existing look (bold/bordered column A, plain data columns).
    $ws.Range("A" + $prevRow + ":P" + $prevRow).Copy() | Out-Null
    $ws.Range("A" + $rowNum + ":P" + $rowNum).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B

    for ($col = 3; $col -le 16; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = 1
    }
}

$excel.CutCopyMode = 0
